$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (Qminus1) values
$ws.Range("B2").Value = 0.1593214268609282
$ws.Range("C2").Value = 0.9530388315012864
$ws.Range("D2").Value = 1.789561637769041
$ws.Range("E2").Value = 1.337744982337456
$ws.Range("F2").Value = 1.343942917851527
$ws.Range("G2").Value = 43

# Update existing row 3 (Q0) values
$ws.Range("B3").Value = 0.2415278289276089
$ws.Range("C3").Value = 1.283805072964911
$ws.Range("D3").Value = 4.178272106320721
$ws.Range("E3").Value = 2.04408221613533
$ws.Range("F3").Value = 2.037554473565015
$ws.Range("G3").Value = 131

# Add new row 4 (Q1) - copy formatting from row 3's label cell so it matches
# the existing bold/centered/bordered header-style used for A2/A3
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A4").Value = "Q1"
$ws.Range("B4").Value = 0.2686780470899982
$ws.Range("C4").Value = 1.393831517495066
$ws.Range("D4").Value = 9.104043377230191
$ws.Range("E4").Value = 3.017290734621076
$ws.Range("F4").Value = 3.029838063148228
$ws.Range("G4").Value = 62
